# Apply the weekly fruit/vegetable update to the "Uva" sheet.
#
# Summary of the change:
#  - Row 117 (Red Globe / $6500-7000-6750 / $/bandeja) is updated in-place to
#    become a "Flame Seedless" / $/caja entry with new date (2022-02-03 / 44595)
#    and new price values.
#  - Three new rows (118, 119, 120) are appended below it:
#       118: Red Globe,        $/caja, date 44595
#       119: Superior Seedless, $/caja, date 44595
#       120: Red Globe,        $/bandeja, date 44335 (the original row 117 data)
#  - The sheet dimension grows from A1:T117 to A1:T120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 120: re-create the original row 117 content (Red Globe / bandeja) ----
$ws.Cells.Item(120, 1).Value = 2
$ws.Cells.Item(120, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(120, 3).Value = "Coquimbo"
$ws.Cells.Item(120, 4).Value = 44335
$ws.Cells.Item(120, 4).NumberFormat = $ws.Cells.Item(117, 4).NumberFormat
$ws.Cells.Item(120, 5).Value = 4
$ws.Cells.Item(120, 6).Value = "Fruta"
$ws.Cells.Item(120, 7).Value = 100109
$ws.Cells.Item(120, 8).Value = "Uva"
$ws.Cells.Item(120, 9).Value = 100109001
$ws.Cells.Item(120, 10).Value = "Uva"
$ws.Cells.Item(120, 11).Value = "Red Globe"
$ws.Cells.Item(120, 12).Value = "Primera"
$ws.Cells.Item(120, 13).Value = 700
$ws.Cells.Item(120, 14).Value = 6500
$ws.Cells.Item(120, 15).Value = 7000
$ws.Cells.Item(120, 16).Value = 6750
$ws.Cells.Item(120, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(120, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(120, 19).Value = 375
$ws.Cells.Item(120, 20).Value = 18

# ---- Row 118: new Red Globe / caja row ----
$ws.Cells.Item(118, 1).Value = 2
$ws.Cells.Item(118, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(118, 3).Value = "Coquimbo"
$ws.Cells.Item(118, 4).Value = 44595
$ws.Cells.Item(118, 4).NumberFormat = $ws.Cells.Item(117, 4).NumberFormat
$ws.Cells.Item(118, 5).Value = 4
$ws.Cells.Item(118, 6).Value = "Fruta"
$ws.Cells.Item(118, 7).Value = 100109
$ws.Cells.Item(118, 8).Value = "Uva"
$ws.Cells.Item(118, 9).Value = 100109001
$ws.Cells.Item(118, 10).Value = "Uva"
$ws.Cells.Item(118, 11).Value = "Red Globe"
$ws.Cells.Item(118, 12).Value = "Primera"
$ws.Cells.Item(118, 13).Value = 400
$ws.Cells.Item(118, 14).Value = 9500
$ws.Cells.Item(118, 15).Value = 10000
$ws.Cells.Item(118, 16).Value = 9750
$ws.Cells.Item(118, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(118, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(118, 19).Value = 542
$ws.Cells.Item(118, 20).Value = 18

# ---- Row 119: new Superior Seedless / caja row ----
$ws.Cells.Item(119, 1).Value = 2
$ws.Cells.Item(119, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(119, 3).Value = "Coquimbo"
$ws.Cells.Item(119, 4).Value = 44595
$ws.Cells.Item(119, 4).NumberFormat = $ws.Cells.Item(117, 4).NumberFormat
$ws.Cells.Item(119, 5).Value = 4
$ws.Cells.Item(119, 6).Value = "Fruta"
$ws.Cells.Item(119, 7).Value = 100109
$ws.Cells.Item(119, 8).Value = "Uva"
$ws.Cells.Item(119, 9).Value = 100109001
$ws.Cells.Item(119, 10).Value = "Uva"
$ws.Cells.Item(119, 11).Value = "Superior Seedless"
$ws.Cells.Item(119, 12).Value = "Primera"
$ws.Cells.Item(119, 13).Value = 400
$ws.Cells.Item(119, 14).Value = 9500
$ws.Cells.Item(119, 15).Value = 10000
$ws.Cells.Item(119, 16).Value = 9750
$ws.Cells.Item(119, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(119, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(119, 19).Value = 542
$ws.Cells.Item(119, 20).Value = 18

# ---- Row 117: update in place (Red Globe -> Flame Seedless, bandeja -> caja) ----
$ws.Cells.Item(117, 4).Value = 44595
$ws.Cells.Item(117, 11).Value = "Flame Seedless"
$ws.Cells.Item(117, 14).Value = 7500
$ws.Cells.Item(117, 15).Value = 8000
$ws.Cells.Item(117, 16).Value = 7750
$ws.Cells.Item(117, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(117, 19).Value = 431
